# Commit: "this is my cart page"
#
# - Rename the worksheet "register " (trailing space) -> "register"
# - Move the active selection from B13 to B23
#
# (The xWindow/yWindow workbook-window position and the revision/document
# GUIDs in the diff are Excel-session/save-time artifacts that aren't part
# of the addressable Excel object model - there is no Range/Worksheet/
# Workbook property that round-trips them - so they are left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the trailing space from the sheet name.
$ws.Name = "register"

# Move the selection to B23 (was B13).
$ws.Range("B23").Select()
